$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 10.97786194618413
$ws.Range("B2").Value = [double]"1.110223024625157e-16"
$ws.Range("C2").Value = 0.00497377543281906
$ws.Range("D2").Value = 0.8429238163007154
$ws.Range("E2").Value = 0.7105205600869624
$ws.Range("A3").Value = 8.381020510234951
$ws.Range("B3").Value = [double]"1.110223024625157e-16"
$ws.Range("C3").Value = 0.005708139113102525
$ws.Range("D3").Value = 0.9673791006814036
$ws.Range("E3").Value = 0.9358223244351614
$ws.Range("A4").Value = 10.47514609081691
$ws.Range("B4").Value = [double]"1.110223024625157e-16"
$ws.Range("C4").Value = 0.005553237262322131
$ws.Range("D4").Value = 0.9411273205245645
$ws.Range("E4").Value = 0.8857206334377464
$ws.Range("A5").Value = 10.09077903288517
$ws.Range("B5").Value = [double]"1.110223024625157e-16"
$ws.Range("C5").Value = 0.006736181214417914
$ws.Range("D5").Value = 1.141605135423672
$ws.Range("E5").Value = 1.303262285225701
$ws.Range("A6").Value = 7.96805397798502
$ws.Range("B6").Value = [double]"1.110223024625157e-16"
$ws.Range("C6").Value = 0.005732176129669167
$ws.Range("D6").Value = 0.9714527413212976
$ws.Range("E6").Value = 0.943720428620664
